# Add Abbyy embedded OCR: trim the extracted-field sheets down to the
# reduced set of fields that the new OCR profile produces.
#
# "Simple Fields" / "Simple Fields - Formatted" sheets keep only:
#   name, vendor-addr, billing-name, billing-addr, invoice-no, date,
#   total, net-amount, items
# (dropping: shipping-addr, po-no, vendor-vat-no, tax, payment-terms,
#  due-date, discount, shipping-charges, payment-addr, currency)
#
# "items" / "items - Formatted" sheets keep only:
#   description, line-amount
# (dropping: quantity, unit-price, item-po-no, line-no, part-no)

$wb = $excel.ActiveWorkbook

# Columns to remove from the two "Simple Fields" sheets, in the order
# shipping-addr(E), po-no(G), vendor-vat-no(H), tax(J), payment-terms(L),
# due-date(N), discount(O), shipping-charges(P), payment-addr(Q), currency(S)
$simpleFieldsColsToDelete = @("S","Q","P","O","N","L","J","H","G","E")

# Columns to remove from the two "items" sheets, in the order
# quantity(B), unit-price(C), item-po-no(E), line-no(F), part-no(G)
$itemsColsToDelete = @("G","F","E","C","B")

foreach ($sheetName in @("Simple Fields", "Simple Fields - Formatted")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($col in $simpleFieldsColsToDelete) {
        $ws.Range("$col`1:$col`2").EntireColumn.Delete()
    }
}

foreach ($sheetName in @("items", "items - Formatted")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($col in $itemsColsToDelete) {
        $ws.Range("$col`1:$col`2").EntireColumn.Delete()
    }
}
